$d = $word.ActiveDocument

function Replace-ExactText($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($rng.Find.Found) {
        # Assign text directly on the matched range so the original
        # apostrophe / quote characters are preserved verbatim (Find's
        # Replace:= argument otherwise "smart-quotes" straight apostrophes).
        $rng.Text = $replaceText
    } else {
        Write-Output "WARNING: text not found: $findText"
    }
}

# le fonds soict rouge. Et lors tu y mectras la  charge de
#   -> le fonds soict rouge. Et lors tu y mectras la susdicte charge de
Replace-ExactText "le fonds soict rouge. Et lors tu y mectras la  charge de" "le fonds soict rouge. Et lors tu y mectras la susdicte charge de"

# tu laisseras consomer le  -> tu laisseras consommer le
Replace-ExactText "tu laisseras consomer le " "tu laisseras consommer le "

# ui est encore meilleur -> ui est encores meilleur
Replace-ExactText "ui est encore meilleur" "ui est encores meilleur"

# sec de telle sorte qu'elles ne sont point pressées
#   -> sec de telle sorte qu'elles ne soient point pressées
Replace-ExactText "sec de telle sorte qu'elles ne sont point pressées" "sec de telle sorte qu'elles ne soient point pressées"

# Il est mieux que parmy le lumignon il y aye un fil de
#   -> Il est mieulx que parmy le lumignon il y aye un fil de
Replace-ExactText "Il est mieux que parmy le lumignon il y aye un fil de " "Il est mieulx que parmy le lumignon il y aye un fil de "
